$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2409
$ws.Range("F6").Value = 59
$ws.Range("F7").Value = 272
$ws.Range("F8").Value = 331
$ws.Range("F9").Value = 2210
$ws.Range("F10").Value = 1156
$ws.Range("F11").Value = 1044
$ws.Range("F12").Value = 848
$ws.Range("F13").Value = 90
$ws.Range("F14").Value = 839
$ws.Range("F15").Value = 1470
$ws.Range("F16").Value = 718
$ws.Range("F17").Value = 1691
$ws.Range("F19").Value = 355
$ws.Range("F20").Value = 65
$ws.Range("F21").Value = 104
$ws.Range("F23").Value = 2615

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 74
$ws.Range("F14").Value = 4
$ws.Range("F21").Value = 4
$ws.Range("F24").Value = 88
$ws.Range("F38").Value = 339
$ws.Range("F43").Value = 75

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2482
$ws.Range("F5").Value = 712
$ws.Range("F6").Value = 2493
$ws.Range("F7").Value = 9533
$ws.Range("F9").Value = 133
$ws.Range("F11").Value = 9
$ws.Range("F12").Value = 356
$ws.Range("F13").Value = 2771
$ws.Range("F14").Value = 358
$ws.Range("F15").Value = 664

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2482
$ws.Range("F3").Value = 712
$ws.Range("F4").Value = 133
$ws.Range("F7").Value = 2409
$ws.Range("F8").Value = 2771
$ws.Range("F9").Value = 358
$ws.Range("F11").Value = 664
$ws.Range("F16").Value = 59
$ws.Range("F17").Value = 272
$ws.Range("F18").Value = 331
$ws.Range("F20").Value = 1044
$ws.Range("F21").Value = 848
$ws.Range("F22").Value = 90
$ws.Range("F23").Value = 839
$ws.Range("F24").Value = 4
$ws.Range("F28").Value = 718
$ws.Range("F31").Value = 1691
$ws.Range("F32").Value = 355
$ws.Range("F33").Value = 88
$ws.Range("F39").Value = 65
$ws.Range("F42").Value = 339
$ws.Range("F43").Value = 2615
$ws.Range("F46").Value = 75
